# Delete the row containing item 20092989 / "AICE ICE SUNDAE 100"
# (originally row 21) so that subsequent rows shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(21).Delete()
